# chore: update Sheets via scheduled runner
# Refreshed market-price-derived figures (currentAveragePrice*, LevePrice*,
# LeveProfit*) for the affected leve rows across the ALC/ARM/BSM/CRP/CUL/
# GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 619.913
$ws.Range("J80").Value = 669.7778
$ws.Range("L80").Value = 2009.3334
$ws.Range("N80").Value = -4005.3334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 619.913
$ws.Range("J83").Value = 669.7778
$ws.Range("L83").Value = 6028.000199999999
$ws.Range("N83").Value = -16012.0002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 1885.4445
$ws.Range("I118").Value = 1000
$ws.Range("J118").Value = 1996.125
$ws.Range("K118").Value = 3000
$ws.Range("L118").Value = 5988.375
$ws.Range("M118").Value = -1343
$ws.Range("N118").Value = -9302.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 847.46344
$ws.Range("J129").Value = 870.8158
$ws.Range("L129").Value = 2612.4474
$ws.Range("N129").Value = -12612.4474

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 741.4375
$ws.Range("I2").Value = 717.7143
$ws.Range("J2").Value = 786.7273
$ws.Range("K2").Value = 717.7143
$ws.Range("L2").Value = 786.7273
$ws.Range("M2").Value = -604.7143
$ws.Range("N2").Value = -1012.7273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1770.46
$ws.Range("I32").Value = 1520.7333
$ws.Range("J32").Value = 4018
$ws.Range("K32").Value = 1520.7333
$ws.Range("L32").Value = 4018
$ws.Range("M32").Value = -1233.7333
$ws.Range("N32").Value = -4592

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 741.4375
$ws.Range("I116").Value = 717.7143
$ws.Range("J116").Value = 786.7273
$ws.Range("K116").Value = 717.7143
$ws.Range("L116").Value = 786.7273
$ws.Range("M116").Value = 1576.2857
$ws.Range("N116").Value = -5374.7273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2479
$ws.Range("J122").Value = 2842.6
$ws.Range("L122").Value = 8527.799999999999
$ws.Range("N122").Value = -13427.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 741.4375
$ws.Range("I3").Value = 717.7143
$ws.Range("J3").Value = 786.7273
$ws.Range("K3").Value = 717.7143
$ws.Range("L3").Value = 786.7273
$ws.Range("M3").Value = -603.7143
$ws.Range("N3").Value = -1014.7273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 12500677
$ws.Range("I94").Value = 19231208
$ws.Range("J94").Value = 1119.8572
$ws.Range("K94").Value = 19231208
$ws.Range("L94").Value = 1119.8572
$ws.Range("M94").Value = -19230757
$ws.Range("N94").Value = -2021.8572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1139.3704
$ws.Range("I107").Value = 922.9167
$ws.Range("K107").Value = 922.9167
$ws.Range("M107").Value = 997.0833

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1165.7391
$ws.Range("I134").Value = 937.6316
$ws.Range("J134").Value = 2249.25
$ws.Range("K134").Value = 2812.8948
$ws.Range("L134").Value = 6747.75
$ws.Range("M134").Value = -277.8948
$ws.Range("N134").Value = -11817.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 3200
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 3200
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 3200
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -3982

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H47").Value = 5000
$ws.Range("J47").Value = 5000
$ws.Range("L47").Value = 5000
$ws.Range("N47").Value = -6132

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 4500
$ws.Range("J48").Value = 4500
$ws.Range("L48").Value = 4500
$ws.Range("N48").Value = -5452

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H49").Value = 3200
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 3200
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 3200
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -3564

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 23557.143
$ws.Range("I51").Value = 20000
$ws.Range("J51").Value = 24150
$ws.Range("K51").Value = 20000
$ws.Range("L51").Value = 24150
$ws.Range("M51").Value = -19264
$ws.Range("N51").Value = -25622

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 23557.143
$ws.Range("I61").Value = 20000
$ws.Range("J61").Value = 24150
$ws.Range("K61").Value = 20000
$ws.Range("L61").Value = 24150
$ws.Range("M61").Value = -19652
$ws.Range("N61").Value = -24846

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1398.875
$ws.Range("I99").Value = 1341.5714
$ws.Range("K99").Value = 1341.5714
$ws.Range("M99").Value = 156.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1398.875
$ws.Range("I126").Value = 1341.5714
$ws.Range("K126").Value = 4024.7142
$ws.Range("M126").Value = -1554.7142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H131").Value = 9102.799999999999
$ws.Range("J131").Value = 16313
$ws.Range("L131").Value = 16313
$ws.Range("N131").Value = -26393

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1379.0968
$ws.Range("I132").Value = 967.0909
$ws.Range("J132").Value = 2386.2222
$ws.Range("K132").Value = 2901.2727
$ws.Range("L132").Value = 7158.6666
$ws.Range("M132").Value = -371.2727
$ws.Range("N132").Value = -12218.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3728.6155
$ws.Range("J39").Value = 3715.6365
$ws.Range("L39").Value = 11146.9095
$ws.Range("N39").Value = -11734.9095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1889.8572
$ws.Range("I122").Value = 727
$ws.Range("J122").Value = 2083.6667
$ws.Range("K122").Value = 6543
$ws.Range("L122").Value = 18753.0003
$ws.Range("M122").Value = -4093
$ws.Range("N122").Value = -23653.0003

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 2584.6667
$ws.Range("I123").Value = 1503.3334
$ws.Range("J123").Value = 2893.6191
$ws.Range("K123").Value = 4510.0002
$ws.Range("L123").Value = 8680.8573
$ws.Range("M123").Value = -2060.0002
$ws.Range("N123").Value = -13580.8573

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 3010
$ws.Range("I125").Value = 1015
$ws.Range("K125").Value = 3045
$ws.Range("M125").Value = 1875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 16129937
$ws.Range("J131").Value = 1063.289
$ws.Range("L131").Value = 3189.867
$ws.Range("N131").Value = -13269.867

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2357.525
$ws.Range("I132").Value = 2096.6072
$ws.Range("J132").Value = 2966.3333
$ws.Range("K132").Value = 6289.821599999999
$ws.Range("L132").Value = 8898.999899999999
$ws.Range("M132").Value = -3759.821599999999
$ws.Range("N132").Value = -13958.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2319.5
$ws.Range("I7").Value = 2284.8572
$ws.Range("J7").Value = 2400.3333
$ws.Range("K7").Value = 2284.8572
$ws.Range("L7").Value = 2400.3333
$ws.Range("M7").Value = -2172.8572
$ws.Range("N7").Value = -2624.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2319.5
$ws.Range("I126").Value = 2284.8572
$ws.Range("J126").Value = 2400.3333
$ws.Range("K126").Value = 6854.571599999999
$ws.Range("L126").Value = 7200.999899999999
$ws.Range("M126").Value = -4384.571599999999
$ws.Range("N126").Value = -12140.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2700.2273
$ws.Range("I132").Value = 2400.1428
$ws.Range("K132").Value = 7200.428400000001
$ws.Range("M132").Value = -4670.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 289.8
$ws.Range("I113").Value = 222.8
$ws.Range("J113").Value = 356.8
$ws.Range("K113").Value = 668.4000000000001
$ws.Range("L113").Value = 1070.4
$ws.Range("M113").Value = 1501.6
$ws.Range("N113").Value = -5410.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 12501386
$ws.Range("I122").Value = 14707273
$ws.Range("J122").Value = 1358.3334
$ws.Range("K122").Value = 44121819
$ws.Range("L122").Value = 4075.0002
$ws.Range("M122").Value = -44119369
$ws.Range("N122").Value = -8975.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 83338984
$ws.Range("I126").Value = 125003600
$ws.Range("K126").Value = 375010800
$ws.Range("M126").Value = -375008330

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4848.125
$ws.Range("I132").Value = 4957.4
$ws.Range("J132").Value = 4666
$ws.Range("K132").Value = 14872.2
$ws.Range("L132").Value = 13998
$ws.Range("M132").Value = -12342.2
$ws.Range("N132").Value = -19058
